# Update Name of Algo
# Apply updated KNN-imputed values to column B for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = 6.127000000000001
    14  = 5.683999999999999
    21  = 9.379000000000001
    23  = 7.886
    25  = 5.67
    26  = 5.701999999999999
    29  = 5.664
    53  = 6.504
    57  = 5.034000000000001
    59  = 4.682
    69  = 5.326
    79  = 5.577000000000001
    83  = 5.824
    91  = 5.601000000000001
    93  = 4.973
    103 = 5.267
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
